# Monthly_Budget_Trunc.xlsx: align the workbook with the example used in the
# paper — rename the "Mortgage or rent" budget line to "Rent" and reword the
# closing question from "Can I afford a fancy dinner?" to
# "Fancy dinner tonight?".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Rent"
$ws.Range("A12").Value = "Fancy dinner tonight?"
